$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text: volume/number + report week dates (new weekly collection) ----
$ws.Range("A8").Value = "Volume 30   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/18/2023  Through  12/24/2023"

# ---- Crime-complaint grid refresh (rows 14-29) ----
$ws.Range("G14").Value = "'0"
$ws.Range("F14").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("H14").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("G15").Value = 2
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -65.217391304347
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 2
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 27.272727272727
$ws.Range("I16").Value = 157
$ws.Range("J16").Value = 185
$ws.Range("K16").Value = -15.135135135135
$ws.Range("L16").Value = 50.961538461538
$ws.Range("M16").Value = -3.680981595092
$ws.Range("N16").Value = -74.42996742671
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = -53.333333333333
$ws.Range("I17").Value = 194
$ws.Range("J17").Value = 211
$ws.Range("K17").Value = -8.056872037914
$ws.Range("L17").Value = 25.974025974026
$ws.Range("M17").Value = 64.406779661017
$ws.Range("N17").Value = -27.340823970037
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("C16").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").Value = 100
$ws.Range("E16").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 125
$ws.Range("I18").Value = 139
$ws.Range("J18").Value = 124
$ws.Range("K18").Value = 12.096774193548
$ws.Range("L18").Value = 49.462365591397
$ws.Range("M18").Value = -16.265060240963
$ws.Range("N18").Value = -87.699115044247
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = 6.666666666666
$ws.Range("F19").Value = 74
$ws.Range("G19").Value = 65
$ws.Range("H19").Value = 13.846153846153
$ws.Range("I19").Value = 686
$ws.Range("J19").Value = 702
$ws.Range("K19").Value = -2.279202279202
$ws.Range("L19").Value = 42.619542619542
$ws.Range("M19").Value = 92.696629213483
$ws.Range("N19").Value = 52.10643015521
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 150
$ws.Range("F20").Value = 29
$ws.Range("G20").Value = 32
$ws.Range("H20").Value = -9.375
$ws.Range("I20").Value = 345
$ws.Range("J20").Value = 279
$ws.Range("K20").Value = 23.655913978494
$ws.Range("L20").Value = 128.476821192053
$ws.Range("M20").Value = 159.398496240602
$ws.Range("N20").Value = -80.330672748004
$ws.Range("C21").Value = 36
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = 33.333333333333
$ws.Range("F21").Value = 133
$ws.Range("G21").Value = 129
$ws.Range("H21").Value = 3.100775193798
$ws.Range("I21").Value = 1534
$ws.Range("J21").Value = 1526
$ws.Range("K21").Value = 0.524246395806
$ws.Range("L21").Value = 53.4
$ws.Range("M21").Value = 61.643835616438
$ws.Range("N21").Value = -63.905882352941
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 1
$ws.Range("C16").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").Value = 100
$ws.Range("E16").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = -20
$ws.Range("I23").Value = 76
$ws.Range("J23").Value = 72
$ws.Range("K23").Value = 5.555555555555
$ws.Range("L23").Value = 72.727272727272
$ws.Range("M23").Value = 72.727272727272
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -8.695652173913
$ws.Range("F24").Value = 92
$ws.Range("G24").Value = 96
$ws.Range("H24").Value = -4.166666666666
$ws.Range("I24").Value = 1140
$ws.Range("J24").Value = 1178
$ws.Range("K24").Value = -3.225806451612
$ws.Range("L24").Value = 27.946127946127
$ws.Range("M24").Value = 43.036386449184
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -45.454545454545
$ws.Range("F25").Value = 31
$ws.Range("G25").Value = 27
$ws.Range("H25").Value = 14.814814814814
$ws.Range("I25").Value = 316
$ws.Range("J25").Value = 328
$ws.Range("K25").Value = -3.658536585365
$ws.Range("L25").Value = 8.965517241379
$ws.Range("M25").Value = -19.38775510204
$ws.Range("C26").Value = "'0"
$ws.Range("D26").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -50
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 36
$ws.Range("K27").Value = -8.333333333333
$ws.Range("L27").Value = 10
$ws.Range("G28").Value = "'0"
$ws.Range("F28").Copy()
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range("H28").Value = "'***.*"
$ws.Range("E28").Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("G29").Value = "'0"
$ws.Range("F29").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("H29").Value = "'***.*"
$ws.Range("E29").Copy()
$ws.Range("H29").PasteSpecial(-4122)

$excel.CutCopyMode = $false
